$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Natalia Galliani"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "gallianinatalia@gmail.com"
$ws.Range("E9").Value = "NataliaGalliani"

$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:gallianinatalia@gmail.com")

$ws.Range("D8").Copy($ws.Range("D9"))
$ws.Range("D9").Value = "gallianinatalia@gmail.com"

$ws.Range("D11").Select()
